$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-07 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-08 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("57÷9=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "22÷5=4, 2", 2) | Out-Null
$d.Content.Find.Execute("51÷2=25, 1", $true, $false, $false, $false, $false, $true, 1, $false, "96÷6=16, 0", 2) | Out-Null
$d.Content.Find.Execute("80÷6=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "50÷4=12, 2", 2) | Out-Null
$d.Content.Find.Execute("36÷7=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "99÷9=11, 0", 2) | Out-Null
$d.Content.Find.Execute("34÷5=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "50÷8=6, 2", 2) | Out-Null
$d.Content.Find.Execute("21÷3=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "26÷5=5, 1", 2) | Out-Null
$d.Content.Find.Execute("57÷8=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "81÷7=11, 4", 2) | Out-Null
$d.Content.Find.Execute("86÷5=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "75÷5=15, 0", 2) | Out-Null
$d.Content.Find.Execute("85÷9=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "35÷7=5, 0", 2) | Out-Null
$d.Content.Find.Execute("51÷9=5, 6", $true, $false, $false, $false, $false, $true, 1, $false, "13÷6=2, 1", 2) | Out-Null
$d.Content.Find.Execute("58÷3=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "66÷4=16, 2", 2) | Out-Null
$d.Content.Find.Execute("23÷9=2, 5", $true, $false, $false, $false, $false, $true, 1, $false, "73÷8=9, 1", 2) | Out-Null
$d.Content.Find.Execute("99÷8=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "54÷6=9, 0", 2) | Out-Null
$d.Content.Find.Execute("78÷4=19, 2", $true, $false, $false, $false, $false, $true, 1, $false, "15÷9=1, 6", 2) | Out-Null
$d.Content.Find.Execute("76÷7=10, 6", $true, $false, $false, $false, $false, $true, 1, $false, "84÷4=21, 0", 2) | Out-Null
$d.Content.Find.Execute("88÷6=14, 4", $true, $false, $false, $false, $false, $true, 1, $false, "17÷8=2, 1", 2) | Out-Null
$d.Content.Find.Execute("75÷4=18, 3", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=5, 0", 2) | Out-Null
$d.Content.Find.Execute("52÷4=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "94÷4=23, 2", 2) | Out-Null
$d.Content.Find.Execute("22÷3=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "54÷6=9, 0", 2) | Out-Null
$d.Content.Find.Execute("90÷8=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "16÷4=4, 0", 2) | Out-Null
$d.Content.Find.Execute("21÷6=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "37÷4=9, 1", 2) | Out-Null
$d.Content.Find.Execute("43÷5=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "66÷4=16, 2", 2) | Out-Null
$d.Content.Find.Execute("75÷9=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "41÷3=13, 2", 2) | Out-Null
$d.Content.Find.Execute("16÷3=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "10÷9=1, 1", 2) | Out-Null
$d.Content.Find.Execute("74÷2=37, 0", $true, $false, $false, $false, $false, $true, 1, $false, "48÷8=6, 0", 2) | Out-Null
